# Updating the files with the Sprint details
$wb = $excel.ActiveWorkbook

# --- Backlog sheet: mark US01 / US02 as Completed (was Planned) ---
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("E10").Value = "Completed"
$backlog.Range("E11").Value = "Completed"

# --- Sprint1 sheet: US01 / US02 rows completed with actual size/time/date ---
$sprint1 = $wb.Worksheets.Item("Sprint1")

# US01 "Date before current dates" row (row 5)
$sprint1.Range("D5").Value = "Completed"
$sprint1.Range("G5").Value = 15
$sprint1.Range("H5").Value = 60
$sprint1.Range("I5").Value = 42646

# US02 "Birth before marriage" row (row 6) - renamed to paired-programming story
$sprint1.Range("B6").Value = "Birth before marriage(Paired Programming)"
$sprint1.Range("C6").Value = "AA/DA"
$sprint1.Range("D6").Value = "Completed"
$sprint1.Range("G6").Value = 15
$sprint1.Range("H6").Value = 60
$sprint1.Range("I6").Value = 42646

# --- Burndown README chart: switch date axis format from m/d/yyyy to m/d/yy ---
$readme = $wb.Worksheets.Item("Burndown README")
$chart = $readme.ChartObjects(1).Chart
$chart.Axes(1).TickLabels.NumberFormat = "m/d/yy"

# --- Active sheet / tab selection moves from Sprint1 to Stories ---
$wb.Worksheets.Item("Stories").Activate()
